# Apply scheduled-runner data refresh to Sheets (auto-generated from OOXML diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 284.2
$ws.Cells.Item(28, 9).Value = 295.35
$ws.Cells.Item(28, 10).Value = 239.6
$ws.Cells.Item(28, 11).Value = 295.35
$ws.Cells.Item(28, 12).Value = 239.6
$ws.Cells.Item(28, 13).Value = 189.65
$ws.Cells.Item(28, 14).Value = -1209.6
$ws.Cells.Item(62, 8).Value = 12732.125
$ws.Cells.Item(62, 9).Value = 18397.77
$ws.Cells.Item(62, 11).Value = 18397.77
$ws.Cells.Item(62, 13).Value = -17773.77
$ws.Cells.Item(65, 8).Value = 12732.125
$ws.Cells.Item(65, 9).Value = 18397.77
$ws.Cells.Item(65, 11).Value = 91988.85000000001
$ws.Cells.Item(65, 13).Value = -88868.85000000001
$ws.Cells.Item(107, 8).Value = 1199.0714
$ws.Cells.Item(107, 9).Value = 2090.5
$ws.Cells.Item(107, 10).Value = 530.5
$ws.Cells.Item(107, 11).Value = 2090.5
$ws.Cells.Item(107, 12).Value = 530.5
$ws.Cells.Item(107, 13).Value = -170.5
$ws.Cells.Item(107, 14).Value = -4370.5
$ws.Cells.Item(116, 8).Value = 4178.5
$ws.Cells.Item(116, 9).Value = 4332.3335
$ws.Cells.Item(116, 10).Value = 4024.6667
$ws.Cells.Item(116, 11).Value = 4332.3335
$ws.Cells.Item(116, 12).Value = 4024.6667
$ws.Cells.Item(116, 13).Value = -890.3334999999997
$ws.Cells.Item(116, 14).Value = -10908.6667
$ws.Cells.Item(125, 8).Value = 620.82355
$ws.Cells.Item(125, 9).Value = 551.8333
$ws.Cells.Item(125, 10).Value = 786.4
$ws.Cells.Item(125, 11).Value = 4966.4997
$ws.Cells.Item(125, 12).Value = 7077.599999999999
$ws.Cells.Item(125, 13).Value = -2506.4997
$ws.Cells.Item(125, 14).Value = -11997.6
$ws.Cells.Item(132, 8).Value = 5123.673
$ws.Cells.Item(132, 9).Value = 1912.3572
$ws.Cells.Item(132, 10).Value = 18611.2
$ws.Cells.Item(132, 11).Value = 5737.071599999999
$ws.Cells.Item(132, 12).Value = 55833.60000000001
$ws.Cells.Item(132, 13).Value = -3207.071599999999
$ws.Cells.Item(132, 14).Value = -60893.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1050
$ws.Cells.Item(45, 9).Value = 900
$ws.Cells.Item(45, 11).Value = 900
$ws.Cells.Item(45, 13).Value = -523
$ws.Cells.Item(61, 8).Value = 671410.2
$ws.Cells.Item(61, 9).Value = 590652.9399999999
$ws.Cells.Item(61, 10).Value = 777015.75
$ws.Cells.Item(61, 11).Value = 590652.9399999999
$ws.Cells.Item(61, 12).Value = 777015.75
$ws.Cells.Item(61, 13).Value = -590440.9399999999
$ws.Cells.Item(61, 14).Value = -777439.75
$ws.Cells.Item(136, 8).Value = 671410.2
$ws.Cells.Item(136, 9).Value = 590652.9399999999
$ws.Cells.Item(136, 10).Value = 777015.75
$ws.Cells.Item(136, 11).Value = 1771958.82
$ws.Cells.Item(136, 12).Value = 2331047.25
$ws.Cells.Item(136, 13).Value = -1769408.82
$ws.Cells.Item(136, 14).Value = -2336147.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 243.8
$ws.Cells.Item(80, 9).Value = 286.57144
$ws.Cells.Item(80, 10).Value = 144
$ws.Cells.Item(80, 11).Value = 286.57144
$ws.Cells.Item(80, 12).Value = 144
$ws.Cells.Item(80, 13).Value = 711.4285600000001
$ws.Cells.Item(80, 14).Value = -2140
$ws.Cells.Item(83, 8).Value = 243.8
$ws.Cells.Item(83, 9).Value = 286.57144
$ws.Cells.Item(83, 10).Value = 144
$ws.Cells.Item(83, 11).Value = 1432.8572
$ws.Cells.Item(83, 12).Value = 720
$ws.Cells.Item(83, 13).Value = 3559.1428
$ws.Cells.Item(83, 14).Value = -10704
$ws.Cells.Item(96, 8).Value = 20613.334
$ws.Cells.Item(105, 8).Value = 1532.6
$ws.Cells.Item(105, 9).Value = 1426.875
$ws.Cells.Item(105, 10).Value = 2660.3333
$ws.Cells.Item(105, 11).Value = 1426.875
$ws.Cells.Item(105, 12).Value = 2660.3333
$ws.Cells.Item(105, 13).Value = 320.125
$ws.Cells.Item(105, 14).Value = -6154.3333
$ws.Cells.Item(107, 8).Value = 1044.7778
$ws.Cells.Item(107, 9).Value = 1044.7778
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 1044.7778
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 875.2221999999999
$ws.Cells.Item(107, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1323.5
$ws.Cells.Item(22, 9).Value = 1522.4445
$ws.Cells.Item(22, 10).Value = 726.6667
$ws.Cells.Item(22, 11).Value = 1522.4445
$ws.Cells.Item(22, 12).Value = 726.6667
$ws.Cells.Item(22, 13).Value = -1172.4445
$ws.Cells.Item(22, 14).Value = -1426.6667
$ws.Cells.Item(31, 8).Value = 1047114.4
$ws.Cells.Item(31, 9).Value = 820.94116
$ws.Cells.Item(31, 10).Value = 2529363.2
$ws.Cells.Item(31, 11).Value = 820.94116
$ws.Cells.Item(31, 12).Value = 2529363.2
$ws.Cells.Item(31, 13).Value = -525.94116
$ws.Cells.Item(31, 14).Value = -2529953.2
$ws.Cells.Item(34, 8).Value = 1047114.4
$ws.Cells.Item(34, 9).Value = 820.94116
$ws.Cells.Item(34, 10).Value = 2529363.2
$ws.Cells.Item(34, 11).Value = 820.94116
$ws.Cells.Item(34, 12).Value = 2529363.2
$ws.Cells.Item(34, 13).Value = -618.94116
$ws.Cells.Item(34, 14).Value = -2529767.2
$ws.Cells.Item(58, 8).Value = 5117.1035
$ws.Cells.Item(58, 9).Value = 6582.3335
$ws.Cells.Item(58, 10).Value = 2719.4546
$ws.Cells.Item(58, 11).Value = 6582.3335
$ws.Cells.Item(58, 12).Value = 2719.4546
$ws.Cells.Item(58, 13).Value = -6379.3335
$ws.Cells.Item(58, 14).Value = -3125.4546
$ws.Cells.Item(107, 8).Value = 308.5946
$ws.Cells.Item(107, 9).Value = 289.43478
$ws.Cells.Item(107, 10).Value = 340.07144
$ws.Cells.Item(107, 11).Value = 289.43478
$ws.Cells.Item(107, 12).Value = 340.07144
$ws.Cells.Item(107, 13).Value = 1630.56522
$ws.Cells.Item(107, 14).Value = -4180.07144
$ws.Cells.Item(136, 8).Value = 5117.1035
$ws.Cells.Item(136, 9).Value = 6582.3335
$ws.Cells.Item(136, 10).Value = 2719.4546
$ws.Cells.Item(136, 11).Value = 19747.0005
$ws.Cells.Item(136, 12).Value = 8158.3638
$ws.Cells.Item(136, 13).Value = -17197.0005
$ws.Cells.Item(136, 14).Value = -13258.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1754.3334
$ws.Cells.Item(14, 9).Value = 1754.3334
$ws.Cells.Item(14, 11).Value = 5263.0002
$ws.Cells.Item(14, 13).Value = -5090.0002
$ws.Cells.Item(92, 8).Value = 335.7143
$ws.Cells.Item(92, 10).Value = 337.5
$ws.Cells.Item(92, 12).Value = 1012.5
$ws.Cells.Item(92, 14).Value = -3508.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(42, 8).Value = 45199
$ws.Cells.Item(42, 10).Value = 45199
$ws.Cells.Item(42, 12).Value = 45199
$ws.Cells.Item(42, 14).Value = -46169
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 10276.333
$ws.Cells.Item(109, 9).Value = 10259
$ws.Cells.Item(109, 11).Value = 10259
$ws.Cells.Item(109, 13).Value = -9219
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 35249.25
$ws.Cells.Item(112, 9).Value = 7000
$ws.Cells.Item(112, 10).Value = 39284.855
$ws.Cells.Item(112, 11).Value = 7000
$ws.Cells.Item(112, 12).Value = 39284.855
$ws.Cells.Item(112, 13).Value = -5892
$ws.Cells.Item(112, 14).Value = -41500.855
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 45199
$ws.Cells.Item(115, 10).Value = 45199
$ws.Cells.Item(115, 12).Value = 45199
$ws.Cells.Item(115, 14).Value = -47549
$ws.Cells.Item(117, 8).Value = 29000
$ws.Cells.Item(117, 10).Value = 29000
$ws.Cells.Item(117, 12).Value = 29000
$ws.Cells.Item(117, 14).Value = -35884
$ws.Cells.Item(119, 8).Value = 49800
$ws.Cells.Item(119, 10).Value = 49800
$ws.Cells.Item(119, 12).Value = 49800
$ws.Cells.Item(119, 14).Value = -59476
$ws.Cells.Item(120, 8).Value = 41500
$ws.Cells.Item(120, 10).Value = 41500
$ws.Cells.Item(120, 12).Value = 41500
$ws.Cells.Item(120, 14).Value = -51176
$ws.Cells.Item(121, 8).Value = 35000
$ws.Cells.Item(121, 10).Value = 35000
$ws.Cells.Item(121, 12).Value = 35000
$ws.Cells.Item(121, 14).Value = -38494
$ws.Cells.Item(122, 8).Value = 1246.7778
$ws.Cells.Item(122, 9).Value = 1246.7778
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3740.3334
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1290.3334
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(123, 8).Value = 45000
$ws.Cells.Item(123, 10).Value = 45000
$ws.Cells.Item(123, 12).Value = 45000
$ws.Cells.Item(123, 14).Value = -49900
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4469115.5
$ws.Cells.Item(132, 9).Value = 9620834
$ws.Cells.Item(132, 10).Value = 4293.3335
$ws.Cells.Item(132, 11).Value = 28862502
$ws.Cells.Item(132, 12).Value = 12880.0005
$ws.Cells.Item(132, 13).Value = -28859972
$ws.Cells.Item(132, 14).Value = -17940.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 350
$ws.Cells.Item(22, 9).Value = 266.66666
$ws.Cells.Item(22, 10).Value = 433.33334
$ws.Cells.Item(22, 11).Value = 266.66666
$ws.Cells.Item(22, 12).Value = 433.33334
$ws.Cells.Item(22, 13).Value = 28.33334000000002
$ws.Cells.Item(22, 14).Value = -1023.33334
$ws.Cells.Item(27, 8).Value = 350
$ws.Cells.Item(27, 9).Value = 266.66666
$ws.Cells.Item(27, 10).Value = 433.33334
$ws.Cells.Item(27, 11).Value = 266.66666
$ws.Cells.Item(27, 12).Value = 433.33334
$ws.Cells.Item(27, 13).Value = -159.66666
$ws.Cells.Item(27, 14).Value = -647.33334
$ws.Cells.Item(46, 8).Value = 873.04
$ws.Cells.Item(46, 9).Value = 712.08105
$ws.Cells.Item(46, 10).Value = 1331.1538
$ws.Cells.Item(46, 11).Value = 712.08105
$ws.Cells.Item(46, 12).Value = 1331.1538
$ws.Cells.Item(46, 13).Value = -524.08105
$ws.Cells.Item(46, 14).Value = -1707.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 680.1539
$ws.Cells.Item(107, 9).Value = 656.2857
$ws.Cells.Item(107, 10).Value = 708
$ws.Cells.Item(107, 11).Value = 1968.8571
$ws.Cells.Item(107, 12).Value = 2124
$ws.Cells.Item(107, 13).Value = -48.85710000000017
$ws.Cells.Item(107, 14).Value = -5964
$ws.Cells.Item(113, 8).Value = 339
$ws.Cells.Item(113, 9).Value = 299.25
$ws.Cells.Item(113, 11).Value = 897.75
$ws.Cells.Item(113, 13).Value = 1272.25
